{"js": "// Locate the sentence that gets simplified (\"and took this opportunity to\n// setup 2FA as well. Opened ...\" -> \"and opened ...\") and gains a trailing\n// \"[repo link]\" reference.\nconst searchText =\n  \"and took this opportunity to setup 2FA as well. Opened a new public repository and made my first commit. \";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  // Replace the whole span with the reworded sentence plus a \"[link]\"\n  // reference to the new public repository, expressed as a raw OOXML\n  // fragment so the hyperlink run/relationship/style come out exactly as\n  // Word would produce them.\n  const ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">\n<w:body>\n<w:p>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">and </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">opened </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>a new public repository and made my first commit.</w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>[</w:t></w:r>\n<w:hyperlink r:id=\"rIdHL1\" w:history=\"1\"><w:r><w:rPr><w:rStyle w:val=\"Hyperlinkki\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>https://n0gear.github.com/lut_frontend</w:t></w:r></w:hyperlink>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>]</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/_rels/document.xml.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rIdHL1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"https://n0gear.github.com/lut_frontend\" TargetMode=\"External\"/>\n</Relationships>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\n  target.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the sentence that is being simplified / getting the repo link appended.\n$searchText = \"and took this opportunity to setup 2FA as well. Opened a new public repository and made my first commit. \"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($searchText)\n\nif ($found) {\n    # Re-wrap the located span in a plain Range so InsertXML replaces its\n    # content instead of appending after it.\n    $target = $d.Range($rng.Start, $rng.End)\n\n    $xml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">\n<w:body>\n<w:p>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">and </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">opened </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>a new public repository and made my first commit.</w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>[</w:t></w:r>\n<w:hyperlink r:id=\"rIdHL1\" w:history=\"1\"><w:r><w:rPr><w:rStyle w:val=\"Hyperlinkki\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>https://n0gear.github.com/lut_frontend</w:t></w:r></w:hyperlink>\n<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>]</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/_rels/document.xml.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n<Relationship Id=\"rIdHL1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"https://n0gear.github.com/lut_frontend\" TargetMode=\"External\"/>\n</Relationships>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n    $target.InsertXML($xml)\n}\n"}
